# ============================================================================
# Edit: attendance_reports/Y2_B2526_GIT_&_Liver_attendance.xlsx
# - Update risk metrics for a batch of students on the "Summary" sheet
#   (percentage/session counts shift because one extra ANATOMY session
#   attendance got recorded for each of them), including a few Status
#   (risk-level) reclassifications.
# - Append the 14 corresponding new attendance-log rows (09/11/2025 ANATOMY
#   session) to the "Attendance" sheet and extend its AutoFilter / used range
#   / hidden _FilterDatabase defined name to match.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet: per-row metric + status updates
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$SummaryRowUpdates = @(
    @{ Row=98; FVal=$null; FSrc=$null; GVal="17.2%"; HVal=17; LVal=5; MVal=2; OVal=2; QVal=1 }
    @{ Row=108; FVal=$null; FSrc=$null; GVal="10.3%"; HVal=19; LVal=3; MVal=4; OVal=1; QVal=1 }
    @{ Row=162; FVal=$null; FSrc=$null; GVal="17.2%"; HVal=17; LVal=5; MVal=2; OVal=2; QVal=1 }
    @{ Row=179; FVal=$null; FSrc=$null; GVal="10.3%"; HVal=19; LVal=3; MVal=4; OVal=2; QVal=1 }
    @{ Row=189; FVal="Low Risk"; FSrc="F12"; GVal="13.8%"; HVal=18; LVal=4; MVal=3; OVal=2; QVal=1 }
    @{ Row=190; FVal="Low Risk"; FSrc="F12"; GVal="13.8%"; HVal=18; LVal=4; MVal=3; OVal=1; QVal=1 }
    @{ Row=215; FVal=$null; FSrc=$null; GVal="17.2%"; HVal=17; LVal=5; MVal=2; OVal=2; QVal=1 }
    @{ Row=219; FVal="No Risk"; FSrc="F3"; GVal="20.7%"; HVal=16; LVal=6; MVal=1; OVal=1; QVal=1 }
    @{ Row=221; FVal="Moderate Risk"; FSrc="F2"; GVal="6.9%"; HVal=20; LVal=2; MVal=5; OVal=1; QVal=1 }
    @{ Row=236; FVal=$null; FSrc=$null; GVal="17.2%"; HVal=17; LVal=5; MVal=2; OVal=2; QVal=1 }
    @{ Row=238; FVal=$null; FSrc=$null; GVal="17.2%"; HVal=17; LVal=5; MVal=2; OVal=2; QVal=1 }
    @{ Row=244; FVal=$null; FSrc=$null; GVal="3.4%"; HVal=21; LVal=1; MVal=6; OVal=1; QVal=1 }
    @{ Row=246; FVal="Low Risk"; FSrc="F12"; GVal="13.8%"; HVal=18; LVal=4; MVal=3; OVal=1; QVal=1 }
    @{ Row=252; FVal=$null; FSrc=$null; GVal="10.3%"; HVal=19; LVal=3; MVal=4; OVal=1; QVal=1 }
)

foreach ($u in $SummaryRowUpdates) {
    $row = $u.Row

    if ($u.FVal -ne $null) {
        # Status text changed risk tier -> copy the fill/font/border format
        # from an existing cell that already carries the target tier's style,
        # then overwrite the text.
        $wsSummary.Range($u.FSrc).Copy()
        $wsSummary.Range("F$row").PasteSpecial(-4122)
        $excel.CutCopyMode = $false
        $wsSummary.Range("F$row").Value = $u.FVal
    }

    # "Percentage" column is stored as literal text (e.g. "17.2%") even
    # though the cell's number format is 0.0% -- flip to text, write, flip
    # back so Excel doesn't reinterpret the string as a numeric percentage.
    $gCell = $wsSummary.Range("G$row")
    $gFmt = $gCell.NumberFormat
    $gCell.NumberFormat = "@"
    $gCell.Value = $u.GVal
    $gCell.NumberFormat = $gFmt

    $wsSummary.Range("H$row").Value = $u.HVal
    $wsSummary.Range("L$row").Value = $u.LVal
    $wsSummary.Range("M$row").Value = $u.MVal
    $wsSummary.Range("O$row").Value = $u.OVal
    $wsSummary.Range("Q$row").Value = $u.QVal
}

# ---------------------------------------------------------------------------
# 2) Attendance sheet: append the 14 new attendance-log rows
# ---------------------------------------------------------------------------
$wsAttendance = $wb.Worksheets.Item("Attendance")

$AttendanceNewRows = @(
    @{ Row=481; Vals=@("222113", "ابرار محمد عبد الله عبد الحميد", "Year 2", "C1", "222113@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:17", "C1") }
    @{ Row=482; Vals=@("221799", "سعاد عبد الماجد احمد عيسى", "Year 2", "C1", "221799@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:21", "C1") }
    @{ Row=483; Vals=@("221779", "اميرة رمضان سلمان ابو جامع", "Year 2", "C1", "221779@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:22", "C1") }
    @{ Row=484; Vals=@("221699", "بيسان محمود محمد عبد الكريم", "Year 2", "C1", "221699@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:28", "C1") }
    @{ Row=485; Vals=@("221606", "راوية الطاهر عبدالله ناصر", "Year 2", "C1", "221606@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:31", "C1") }
    @{ Row=486; Vals=@("221307", "دعاء عاصم على العوض", "Year 2", "C1", "221307@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:34", "C1") }
    @{ Row=487; Vals=@("222028", "هاجر عبد الحفيظ سيد صالح", "Year 2", "C1", "222028@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:22", "C1") }
    @{ Row=488; Vals=@("221022", "شهد عبدالعظيم فرج بابكر", "Year 2", "C1", "221022@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:24", "C1") }
    @{ Row=489; Vals=@("221944", "ضحى عمر سيف الدين محمد", "Year 2", "C1", "221944@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:28", "C1") }
    @{ Row=490; Vals=@("222056", "الغالى ادم عيسى رحيل", "Year 2", "C1", "222056@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:31", "C1") }
    @{ Row=491; Vals=@("222063", "محمد مصطفى حامد التوم", "Year 2", "C1", "222063@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:34", "C1") }
    @{ Row=492; Vals=@("221949", "محمد لطفى الزبير البشير", "Year 2", "C1", "221949@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:17", "C1") }
    @{ Row=493; Vals=@("222026", "عثمان موسى محمد ادم", "Year 2", "C1", "222026@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:28", "C1") }
    @{ Row=494; Vals=@("221930", "احمد محمود عبد الباقى محمد", "Year 2", "C1", "221930@med.asu.edu.eg", "ANATOMY", "2", "ANATOMY", "09/11/2025", "14:23:34", "C1") }
)

$firstNewRow = 481
$lastNewRow = 494

# Columns A (Student ID), G (Session) and I (Date) contain values that look
# numeric/date-like ("222113", "2", "09/11/2025"); force them to be stored
# as literal text (matching every other row in this sheet) instead of
# letting Excel auto-convert them to a number / a mm-dd date serial.
$wsAttendance.Range("A$firstNewRow`:A$lastNewRow").NumberFormat = "@"
$wsAttendance.Range("G$firstNewRow`:G$lastNewRow").NumberFormat = "@"
$wsAttendance.Range("I$firstNewRow`:I$lastNewRow").NumberFormat = "@"

foreach ($r in $AttendanceNewRows) {
    $row = $r.Row
    $vals = $r.Vals
    for ($col = 1; $col -le $vals.Length; $col++) {
        $wsAttendance.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}

# Extend the AutoFilter range to cover the newly-appended rows.
$wsAttendance.AutoFilterMode = $false
$wsAttendance.Range("A1:K$lastNewRow").AutoFilter()

# Keep the workbook-level hidden _xlnm._FilterDatabase name for the
# Attendance sheet in sync with the new AutoFilter extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$$lastNewRow"
    }
}
